$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.104.10"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.788.88"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'226.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("D6").Value = "'0.545"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "'32.25"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("D9").Value = "'0.294"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.31%  "
$ws.Range("D10").Value = "'0.0687"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.71%  "
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").Value = "2.047.75"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "'11.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.88%  "
$ws.Range("D14").Value = "1.784.87"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "34.086.26"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "'0.620"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.50%  "
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").Value = "'244.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").Value = "0.0₃0776"
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("D21").Value = "'11.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.15%  "
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Value = "'2.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.16%  "
$ws.Range("D25").Value = "'159.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("D26").Value = "'7.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'16.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  +0.97%  "
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("D31").Value = "'0.0518"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("E32").Value = "  -0.32%  "
$ws.Range("D33").Value = "'3.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.54%  "
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("D35").Value = "1.406.12"
$ws.Range("D36").Value = "'0.654"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("E37").Value = "  -0.50%  "
$ws.Range("D38").Value = "'0.0188"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.28%  "
$ws.Range("D39").Value = "'2.32"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.83%  "
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("D41").Value = "'79.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.22%  "
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").Value = "'13.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.63%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'6.04"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.74%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0138"
$ws.Range("E46").Value = "  -5.34%  "
$ws.Range("E47").Value = "  +2.40%  "
$ws.Range("D48").Value = "'1.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("D49").Value = "'106.87"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.34%  "
$ws.Range("D50").Value = "1.948.96"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("E51").Value = "  +0.10%  "
